$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "303.16"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2.07%"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "32.18"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.32%"
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.122"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.90%"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07851"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.81%"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.251"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-9.29%"
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.817"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.54%"
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.797"
# Row 9
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.38%"
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1761"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.49%"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07815"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "6.28%"
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08814"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.29%"
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03054"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.80%"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1000"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.02%"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001508"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.11%"
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005999"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.61%"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.464"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.55%"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.250"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.32%"
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1347"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.78%"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.268"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.51%"
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1797"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "11.14%"
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04583"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.49%"
# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.76%"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004497"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.73%"
# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "4.28%"
# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-1.30%"
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01791"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.16%"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04729"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "5.73%"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007238"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.78%"
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1371"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "2.13%"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002126"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-3.88%"
# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "13.03%"
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006205"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-5.46%"
# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.20%"
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.003205"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-38.80%"
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.117"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "36.10%"
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002104"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.20%"
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002004"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.20%"
